$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.697.72'
$ws.Range("E2").Value = '  -0.90%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.527.06'
$ws.Range("E3").Value = '  -3.14%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '203.43'
$ws.Range("E5").Value = '  +2.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '556.70'
$ws.Range("E6").Value = '  -3.87%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.514.02'
$ws.Range("E7").Value = '  -3.32%  '
$ws.Range("E8").Value = '  -1.79%  '
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '64.56'
$ws.Range("E10").Value = '  +13.82%  '
$ws.Range("E11").Value = '  -2.95%  '
$ws.Range("E12").Value = '  -6.52%  '
$ws.Range("E13").Value = '  -7.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.93'
$ws.Range("E14").Value = '  -1.87%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.078.63'
$ws.Range("E15").Value = '  -3.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.518.15'
$ws.Range("E16").Value = '  -3.31%  '
$ws.Range("E17").Value = '  -1.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.56'
$ws.Range("E18").Value = '  -0.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '67.408.87'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.87'
$ws.Range("E20").Value = '  -5.80%  '
$ws.Range("E21").Value = '  -5.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '393.89'
$ws.Range("E22").Value = '  -2.44%  '
$ws.Range("E23").Value = '  -7.37%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.02'
$ws.Range("E24").Value = '  -5.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.42'
$ws.Range("E25").Value = '  -3.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.92'
$ws.Range("E26").Value = '  +1.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.85'
$ws.Range("E27").Value = '  -3.94%  '
$ws.Range("E28").Value = '  -3.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.91'
$ws.Range("E29").Value = '  -3.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '715.57'
$ws.Range("E30").Value = '  +4.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '31.18'
$ws.Range("E31").Value = '  -1.96%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.13'
$ws.Range("E32").Value = '  -13.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.79'
$ws.Range("E33").Value = '  -3.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '64.10'
$ws.Range("E34").Value = '  -0.94%  '
$ws.Range("E35").Value = '  -5.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '38.78'
$ws.Range("E36").Value = '  -9.44%  '
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.400'
$ws.Range("E38").Value = '  -6.70%  '
$ws.Range("E39").Value = '  -4.31%  '
$ws.Range("E40").Value = '  -4.91%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.997'
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.063.22'
$ws.Range("E42").Value = '  -4.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₃0690'
$ws.Range("E43").Value = '  -12.71%  '
$ws.Range("E44").Value = '  -11.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.76'
$ws.Range("E45").Value = '  +5.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.76'
$ws.Range("E46").Value = '  -9.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0409'
$ws.Range("E47").Value = '  -3.04%  '
$ws.Range("E48").Value = '  -3.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '138.41'
$ws.Range("E49").Value = '  -2.39%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.30'
$ws.Range("E50").Value = '  -7.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.86'
$ws.Range("E51").Value = '  -7.95%  '
